{"js": "// Update the title date line and the 20x5 grid of addition/subtraction\n// expressions to the new values from the target revision.\n\nconst newTableValues = [\n  [\"41+0=\", \"20+13=\", \"54+19=\", \"2+91=\", \"16+31=\"],\n  [\"61-9=\", \"47-23=\", \"37+62=\", \"68+4=\", \"11+59=\"],\n  [\"67-22=\", \"90-62=\", \"25+53=\", \"39-30=\", \"52+27=\"],\n  [\"9+11=\", \"19-15=\", \"99-89=\", \"29+20=\", \"8+81=\"],\n  [\"22+21=\", \"41+56=\", \"86-51=\", \"6+82=\", \"13+81=\"],\n  [\"96-23=\", \"77-73=\", \"64-56=\", \"83+16=\", \"76-31=\"],\n  [\"26+59=\", \"77-11=\", \"7+1=\", \"79-14=\", \"37-22=\"],\n  [\"60+5=\", \"33+19=\", \"45-27=\", \"24+64=\", \"29+55=\"],\n  [\"96-56=\", \"71+3=\", \"14+20=\", \"96-27=\", \"44+39=\"],\n  [\"84-81=\", \"48+5=\", \"39+32=\", \"87-34=\", \"54-28=\"],\n  [\"40-6=\", \"77-64=\", \"1+0=\", \"5+37=\", \"62-1=\"],\n  [\"20+47=\", \"5+31=\", \"2+6=\", \"76-52=\", \"5+66=\"],\n  [\"75-71=\", \"15-5=\", \"11+88=\", \"24+68=\", \"71-68=\"],\n  [\"82-34=\", \"73-45=\", \"96-68=\", \"39-30=\", \"5-2=\"],\n  [\"64+34=\", \"90-8=\", \"53-15=\", \"66-28=\", \"73-38=\"],\n  [\"77-47=\", \"34-10=\", \"2+55=\", \"23+46=\", \"60-18=\"],\n  [\"1+56=\", \"89-6=\", \"56-54=\", \"52-5=\", \"18+24=\"],\n  [\"57-2=\", \"98-47=\", \"35+61=\", \"74-42=\", \"97-19=\"],\n  [\"84-44=\", \"49-37=\", \"36-16=\", \"8+58=\", \"56-20=\"],\n  [\"57-23=\", \"94-8=\", \"87-79=\", \"6+78=\", \"61-25=\"],\n];\n\n// 1) Update the title paragraph (date line) above the table.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\nif (titlePara.text.trim() === \"2024-01-14 Sunday\") {\n  titlePara.insertText(\"2024-01-15 Monday\", Word.InsertLocation.replace);\n}\n\n// 2) Update every cell value in the (first) table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount\");\nawait context.sync();\n\ntable.values = newTableValues;\nawait context.sync();\n", "ps1": "# Update the title date line and the 20x5 grid of addition/subtraction\n# expressions to the new values from the target revision.\n\n$d = $word.ActiveDocument\n\n# 1) Update the title paragraph (date line) above the table.\n$titlePara = $d.Paragraphs.Item(1)\nif ($titlePara.Range.Text.TrimEnd([char]13) -eq \"2024-01-14 Sunday\") {\n    $titlePara.Range.Text = \"2024-01-15 Monday\"\n}\n\n# 2) Update every cell value in the (first) table, row by row.\n$newTableValues = @(\n    , @(\"41+0=\", \"20+13=\", \"54+19=\", \"2+91=\", \"16+31=\")\n    , @(\"61-9=\", \"47-23=\", \"37+62=\", \"68+4=\", \"11+59=\")\n    , @(\"67-22=\", \"90-62=\", \"25+53=\", \"39-30=\", \"52+27=\")\n    , @(\"9+11=\", \"19-15=\", \"99-89=\", \"29+20=\", \"8+81=\")\n    , @(\"22+21=\", \"41+56=\", \"86-51=\", \"6+82=\", \"13+81=\")\n    , @(\"96-23=\", \"77-73=\", \"64-56=\", \"83+16=\", \"76-31=\")\n    , @(\"26+59=\", \"77-11=\", \"7+1=\", \"79-14=\", \"37-22=\")\n    , @(\"60+5=\", \"33+19=\", \"45-27=\", \"24+64=\", \"29+55=\")\n    , @(\"96-56=\", \"71+3=\", \"14+20=\", \"96-27=\", \"44+39=\")\n    , @(\"84-81=\", \"48+5=\", \"39+32=\", \"87-34=\", \"54-28=\")\n    , @(\"40-6=\", \"77-64=\", \"1+0=\", \"5+37=\", \"62-1=\")\n    , @(\"20+47=\", \"5+31=\", \"2+6=\", \"76-52=\", \"5+66=\")\n    , @(\"75-71=\", \"15-5=\", \"11+88=\", \"24+68=\", \"71-68=\")\n    , @(\"82-34=\", \"73-45=\", \"96-68=\", \"39-30=\", \"5-2=\")\n    , @(\"64+34=\", \"90-8=\", \"53-15=\", \"66-28=\", \"73-38=\")\n    , @(\"77-47=\", \"34-10=\", \"2+55=\", \"23+46=\", \"60-18=\")\n    , @(\"1+56=\", \"89-6=\", \"56-54=\", \"52-5=\", \"18+24=\")\n    , @(\"57-2=\", \"98-47=\", \"35+61=\", \"74-42=\", \"97-19=\")\n    , @(\"84-44=\", \"49-37=\", \"36-16=\", \"8+58=\", \"56-20=\")\n    , @(\"57-23=\", \"94-8=\", \"87-79=\", \"6+78=\", \"61-25=\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowValues = $newTableValues[$r - 1]\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
